# Applies the "Alpha" absence-type addition to the Permission Request template.
#  - Example sheet: a few sample rows' "Type Izin" (col G) values are updated
#    so the sample data showcases the new option.
#  - Both sheets: the "Type Izin" (col G) dropdown list gains an "Alpha"
#    entry and its applicable range is widened.

$wb = $excel.ActiveWorkbook

$wsPermission = $wb.Worksheets.Item("Permission")
$wsExample    = $wb.Worksheets.Item("Example")

# --- Example sheet: update sample data to use the new "Alpha" option -------
$wsExample.Range("G2").Value = "Izin Terlambat"
$wsExample.Range("G3").Value = "Alpha"
$wsExample.Range("G4").Value = "Alpha"

# --- Permission sheet: widen "Type Izin" dropdown & add "Alpha" ------------
$wsPermission.Range("G2:G44").Validation.Delete()
$wsPermission.Range("G2:G217").Validation.Add(3, 1, 1, '"Izin Terlambat, Izin Tidak Masuk Kerja, Izin Pulang Cepat, Izin Keluar Kantor, Alpha"')

# --- Example sheet: widen "Type Izin" dropdown & add "Alpha" ---------------
$wsExample.Range("G2:G6").Validation.Delete()
$wsExample.Range("G2:G18").Validation.Add(3, 1, 1, '"Izin Terlambat, Izin Tidak Masuk Kerja, Izin Pulang Cepat, Izin Keluar Kantor, Alpha"')

# --- Selection / active sheet, matching the author's final UI state --------
$wsPermission.Range("G2").Select() | Out-Null
$wsExample.Range("F5").Select() | Out-Null
$wsExample.Activate() | Out-Null
